$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property; set its value to the literal text "false"
# (not the boolean FALSE). Typing "false" directly gets auto-coerced to a
# boolean by Excel's input parser, so we enter it with a leading apostrophe
# to force text, then copy just the formatting from the sibling cell (A7)
# back onto B7 so the quote-prefix doesn't leave a stray number format
# behind.
$ws.Range("B7").Value = "'false"
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 8 = "Date" property; update the generated timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
